$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2..26 (Actual Consumption MW, Timestamp)
$data = @(
    @(5279, 45786),
    @(5265, 45786.01041666666),
    @(5297, 45786.02083333334),
    @(5226, 45786.03125),
    @(5183, 45786.04166666666),
    @(5172, 45786.05208333334),
    @(5168, 45786.0625),
    @(5135, 45786.07291666666),
    @(5061, 45786.08333333334),
    @(5097, 45786.09375),
    @(5063, 45786.10416666666),
    @(5097, 45786.11458333334),
    @(5124, 45786.125),
    @(5111, 45786.13541666666),
    @(5157, 45786.14583333334),
    @(5181, 45786.15625),
    @(5240, 45786.16666666666),
    @(5239, 45786.17708333334),
    @(5239, 45786.1875),
    @(5255, 45786.19791666666),
    @(5307, 45786.20833333334),
    @(5417, 45786.21875),
    @(5522, 45786.22916666666),
    @(5657, 45786.23958333334),
    @(5807, 45786.25)
)

# Delete rows 27 through 39 (old extra data beyond the new 26-row extent)
$ws.Range("A27:B39").EntireRow.Delete() | Out-Null

# Write the updated values into rows 2..26
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
